# Figure 4f: add the file-location text after the trailing "Excel file:"
# bullet in the Lipolysis (Figure 4) section - i.e. the short bullet that
# currently contains only "Excel file:" (the one that also carries the
# document's "_GoBack" bookmark).

$d = $word.ActiveDocument

# --- locate the target paragraph (content-based, not index-based) -------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Excel file:") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Excel file:' paragraph to edit"
}

# --- the "_GoBack" bookmark sits right where we need to insert text -----
# InsertXML on a range that spans an existing "_GoBack" bookmark makes the
# host re-wrap that bookmark around the whole replaced range (mirroring how
# Word re-anchors _GoBack on an edit). Removing it first and leaving it out
# of the replacement markup below lets us control exactly where the
# (recreated) bookmark ends up - sandwiched between "qpcr" and the closing
# spell-check marker, matching the authored edit.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$full = $d.Range($target.Range.Start, $target.Range.End)

$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Excel file:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Harvey&gt;Mouse work&gt; Dexamethasone treatment&gt; cohort A</w:t></w:r><w:r><w:t>&gt;lipolysis folder&gt;</w:t></w:r><w:r><w:t xml:space="preserve"> 2016-08-25 152316 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dex</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cohort A lipolysis </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Iwat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>qpcr</w:t></w:r><w:bookmarkStart w:id="3" w:name="_GoBack"/><w:bookmarkEnd w:id="3"/><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$full.InsertXML($xmlFrag)
